$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B62 to be a numeric value (3) instead of text "3"
$ws.Cells.Item(62, 2).Value = 3

# Add new row 63 with the new annotation data
$ws.Cells.Item(63, 1).Value = "Ying Tang"

$ws.Cells.Item(63, 2).NumberFormat = "@"
$ws.Cells.Item(63, 2).Value = "3"
$ws.Cells.Item(63, 2).Style = "Normal"

$ws.Cells.Item(63, 3).Value = "无"
$ws.Cells.Item(63, 4).Value = "FBK"
$ws.Cells.Item(63, 5).Value = "MET"
$ws.Cells.Item(63, 6).Value = "c32ea772-309c-4b62-b4b3-b2a94581f569"
$ws.Cells.Item(63, 7).Value = "S1Euwz-Rb_annotated.xlsx"
$ws.Cells.Item(63, 8).Value = "- For the gating mechanism of the writing unit, we have performed additional experiments showing that untied gate values for each entry of the state vector perform better than having one shared potentially-interpretable gate for the whole state and so have changed the description of that subsection accordingly."
